# Modification des scripts deploiements
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the UUID values held in B2 and B3 for new ones
$ws.Range("B2").Value = "643a506c-55fd-460a-b8aa-6ac6c470b8ea"
$ws.Range("B3").Value = "4e37b55e-4df0-48d4-aefe-75d9d05ed21c"

# Update the active selection to B2 (was B3)
$ws.Range("B2").Select()
